# TC03_C3DC_phs002504_Race-AmerIndAlaskNat.xlsx
# - Fix the "Treatment Agent" query in the Treatment tab (row 5 / cell B5):
#   drop the redundant CONCAT() wrapper around REPLACE(...).
# - Re-apply formatting on B5 (mirrors the font bump Excel recorded when the
#   cell was last edited/re-saved) and leave B5 as the active selection,
#   matching the state the workbook was left in.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cell = $ws.Range("B5")

$oldText = $cell.Value2
$newText = $oldText.Replace( `
    "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent""", `
    "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent""")
$cell.Value2 = $newText

# Touch the font so the cell gets its own (new) style entry, same as the
# authoring session recorded for this cell.
$cell.Font.ThemeColor = 1
$cell.Font.Size = 12
$cell.WrapText = $true

# Leave the selection on B5 (row 5 is scrolled into view).
$cell.Select()

Write-Output "Updated Treatment Agent formula on B5 and refreshed its selection/format."
